# Update "Pais" (countries) COVID data sheet per the 4 Abr 2020 01:20 refresh.
# - Updates the "last updated" timestamp text.
# - Refreshes totals for a handful of countries whose stats changed between
#   the 00:50 and 01:20 snapshots (several of which also shifted row position
#   as the underlying source list got re-sorted / had entries inserted).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 4 de Abril de 2020 a las 01:20"
$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 276037
$ws.Range("C4").Value = 31160
$ws.Range("D4").Value = 12268
$ws.Range("E4").Value = 256384
$ws.Range("F4").Value = 5787
$ws.Range("G4").Value = 1314
$ws.Range("H4").Value = 7385
$ws.Range("A20").Value = "Brasil"
$ws.Range("B20").Value = 9194
$ws.Range("C20").Value = 1150
$ws.Range("D20").Value = 127
$ws.Range("E20").Value = 8704
$ws.Range("F20").Value = 296
$ws.Range("G20").Value = 39
$ws.Range("H20").Value = 363
$ws.Range("A23").Value = "Australia"
$ws.Range("B23").Value = 5454
$ws.Range("C23").Value = 140
$ws.Range("D23").Value = 585
$ws.Range("E23").Value = 4841
$ws.Range("F23").Value = 85
$ws.Range("G23").Value = 3
$ws.Range("H23").Value = 28
$ws.Range("A24").Value = "Noruega"
$ws.Range("B24").Value = 5370
$ws.Range("C24").Value = 223
$ws.Range("D24").Value = 32
$ws.Range("E24").Value = 5279
$ws.Range("F24").Value = 96
$ws.Range("G24").Value = 9
$ws.Range("H24").Value = 59
$ws.Range("A26").Value = "Chequia"
$ws.Range("B26").Value = 4190
$ws.Range("C26").Value = 332
$ws.Range("D26").Value = 72
$ws.Range("E26").Value = 4065
$ws.Range("F26").Value = 77
$ws.Range("G26").Value = 9
$ws.Range("H26").Value = 53
$ws.Range("A27").Value = "Rusia"
$ws.Range("B27").Value = 4149
$ws.Range("C27").Value = 601
$ws.Range("D27").Value = 281
$ws.Range("E27").Value = 3834
$ws.Range("F27").Value = 8
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = 34
$ws.Range("A96").Value = "San Marino"
$ws.Range("B96").Value = 251
$ws.Range("C96").Value = 6
$ws.Range("D96").Value = 26
$ws.Range("E96").Value = 193
$ws.Range("F96").Value = 13
$ws.Range("G96").Value = 2
$ws.Range("H96").Value = 32
$ws.Range("A112").Value = "Consejo Danes para los Refugiados"
$ws.Range("B112").Value = 148
$ws.Range("C112").Value = 14
$ws.Range("D112").Value = 3
$ws.Range("E112").Value = 129
$ws.Range("F112").Value = 0
$ws.Range("G112").Value = 3
$ws.Range("H112").Value = 16
$ws.Range("A113").Value = "Martinica"
$ws.Range("B113").Value = 143
$ws.Range("C113").Value = 5
$ws.Range("D113").Value = 27
$ws.Range("E113").Value = 113
$ws.Range("F113").Value = 18
$ws.Range("G113").Value = 0
$ws.Range("H113").Value = 3
$ws.Range("A134").Value = "Jamaica"
$ws.Range("B134").Value = 53
$ws.Range("C134").Value = 6
$ws.Range("D134").Value = 7
$ws.Range("E134").Value = 43
$ws.Range("F134").Value = 0
$ws.Range("G134").Value = 0
$ws.Range("H134").Value = 3
$ws.Range("A135").Value = "Barbados"
$ws.Range("B135").Value = 51
$ws.Range("C135").Value = 5
$ws.Range("D135").Value = 0
$ws.Range("E135").Value = 51
$ws.Range("F135").Value = 0
$ws.Range("G135").Value = 0
$ws.Range("H135").Value = 0
$ws.Range("A136").Value = "Guatemala"
$ws.Range("B136").Value = 50
$ws.Range("C136").Value = 3
$ws.Range("D136").Value = 12
$ws.Range("E136").Value = 37
$ws.Range("F136").Value = 1
$ws.Range("G136").Value = 0
$ws.Range("H136").Value = 1
$ws.Range("A137").Value = "Republica de Yibuti"
$ws.Range("B137").Value = 49
$ws.Range("C137").Value = 9
$ws.Range("D137").Value = 8
$ws.Range("E137").Value = 41
$ws.Range("F137").Value = 0
$ws.Range("G137").Value = 0
$ws.Range("H137").Value = 0
$ws.Range("A138").Value = "Uganda"
$ws.Range("B138").Value = 48
$ws.Range("C138").Value = 3
$ws.Range("D138").Value = 0
$ws.Range("E138").Value = 48
$ws.Range("F138").Value = 0
$ws.Range("G138").Value = 0
$ws.Range("H138").Value = 0
$ws.Range("A139").Value = "El Salvador"
$ws.Range("B139").Value = 46
$ws.Range("C139").Value = 5
$ws.Range("D139").Value = 0
$ws.Range("E139").Value = 43
$ws.Range("F139").Value = 4
$ws.Range("G139").Value = 1
$ws.Range("H139").Value = 3
$ws.Range("A143").Value = "Mali"
$ws.Range("B143").Value = 39
$ws.Range("C143").Value = 3
$ws.Range("D143").Value = 0
$ws.Range("E143").Value = 36
$ws.Range("F143").Value = 0
$ws.Range("G143").Value = 0
$ws.Range("H143").Value = 3
$ws.Range("A144").Value = "Zambia"
$ws.Range("B144").Value = 39
$ws.Range("C144").Value = 0
$ws.Range("D144").Value = 2
$ws.Range("E144").Value = 36
$ws.Range("F144").Value = 0
$ws.Range("G144").Value = 0
$ws.Range("H144").Value = 1
$ws.Range("A145").Value = "Puerto Rico"
$ws.Range("B145").Value = 39
$ws.Range("C145").Value = 0
$ws.Range("D145").Value = 1
$ws.Range("E145").Value = 36
$ws.Range("F145").Value = 0
$ws.Range("G145").Value = 0
$ws.Range("H145").Value = 2
$ws.Range("A150").Value = "Bahamas"
$ws.Range("B150").Value = 24
$ws.Range("C150").Value = 0
$ws.Range("D150").Value = 0
$ws.Range("E150").Value = 21
$ws.Range("F150").Value = 1
$ws.Range("G150").Value = 2
$ws.Range("H150").Value = 3
$ws.Range("A160").Value = "Nueva Caledonia"
$ws.Range("B160").Value = 18
$ws.Range("C160").Value = 0
$ws.Range("D160").Value = 1
$ws.Range("E160").Value = 17
$ws.Range("F160").Value = 0
$ws.Range("G160").Value = 0
$ws.Range("H160").Value = 0
$ws.Range("A161").Value = "Haiti"
$ws.Range("B161").Value = 18
$ws.Range("C161").Value = 2
$ws.Range("D161").Value = 1
$ws.Range("E161").Value = 17
$ws.Range("F161").Value = 0
$ws.Range("G161").Value = 0
$ws.Range("H161").Value = 0
$ws.Range("A163").Value = "Libia"
$ws.Range("B163").Value = 17
$ws.Range("C163").Value = 6
$ws.Range("D163").Value = 0
$ws.Range("E163").Value = 16
$ws.Range("F163").Value = 0
$ws.Range("G163").Value = 0
$ws.Range("H163").Value = 1
$ws.Range("A164").Value = "Guinea Ecuatorial"
$ws.Range("B164").Value = 16
$ws.Range("C164").Value = 1
$ws.Range("D164").Value = 1
$ws.Range("E164").Value = 15
$ws.Range("F164").Value = 0
$ws.Range("G164").Value = 0
$ws.Range("H164").Value = 0
$ws.Range("A166").Value = "Siria"
$ws.Range("B166").Value = 16
$ws.Range("C166").Value = 0
$ws.Range("D166").Value = 0
$ws.Range("E166").Value = 14
$ws.Range("F166").Value = 0
$ws.Range("G166").Value = 0
$ws.Range("H166").Value = 2
$ws.Range("A168").Value = "Guinea-Bisau"
$ws.Range("B168").Value = 15
$ws.Range("C168").Value = 6
$ws.Range("D168").Value = 0
$ws.Range("E168").Value = 15
$ws.Range("F168").Value = 0
$ws.Range("G168").Value = 0
$ws.Range("H168").Value = 0
$ws.Range("A169").Value = "Mongolia"
$ws.Range("B169").Value = 14
$ws.Range("C169").Value = 0
$ws.Range("D169").Value = 2
$ws.Range("E169").Value = 12
$ws.Range("F169").Value = 0
$ws.Range("G169").Value = 0
$ws.Range("H169").Value = 0
$ws.Range("A170").Value = "Namibia"
$ws.Range("B170").Value = 14
$ws.Range("C170").Value = 0
$ws.Range("D170").Value = 3
$ws.Range("E170").Value = 11
$ws.Range("F170").Value = 0
$ws.Range("G170").Value = 0
$ws.Range("H170").Value = 0
$ws.Range("A171").Value = "Santa Lucia"
$ws.Range("B171").Value = 13
$ws.Range("C171").Value = 0
$ws.Range("D171").Value = 1
$ws.Range("E171").Value = 12
$ws.Range("F171").Value = 0
$ws.Range("G171").Value = 0
$ws.Range("H171").Value = 0
$ws.Range("A173").Value = "Dominica"
$ws.Range("B173").Value = 12
$ws.Range("C173").Value = 0
$ws.Range("D173").Value = 0
$ws.Range("E173").Value = 12
$ws.Range("F173").Value = 0
$ws.Range("G173").Value = 1
$ws.Range("H173").Value = 1
$ws.Range("A175").Value = "Mozambique"
$ws.Range("B175").Value = 10
$ws.Range("C175").Value = 0
$ws.Range("D175").Value = 0
$ws.Range("E175").Value = 10
$ws.Range("F175").Value = 0
$ws.Range("G175").Value = 0
$ws.Range("H175").Value = 0
$ws.Range("A176").Value = "Laos"
$ws.Range("B176").Value = 10
$ws.Range("C176").Value = 0
$ws.Range("D176").Value = 0
$ws.Range("E176").Value = 10
$ws.Range("F176").Value = 0
$ws.Range("G176").Value = 0
$ws.Range("H176").Value = 0
